# Natmi following Dr Hou advice
# Rebuild the LR-pairs result table for F2-F2rl1: two new sending clusters
# (M1, M2) are added and every sending/target cluster combination's
# statistics are refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
  @("ECs", "F2", "F2rl1", "FAPs", 3, 1, 0.9442576666666667, 2.832773, 0.4123198466914758, 0.4405367670065795, 3, 1, 2.662448666666667, 7.987346000000001, 0.9268761553538212, 0.9500328399529846, 2.514037565606444, 22.626338090458, 0.3821694342774721, 0.418524395862967),
  @("ECs", "F2", "F2rl1", "sCs", 3, 1, 0.9442576666666667, 2.832773, 0.4123198466914758, 0.4405367670065795, 2, 1, 0.210048, 0.420096, 0.07312384464617888, 0.04996716004701549, 0.198339434368, 1.190036606208, 0.03015041241400377, 0.02201237114361253),
  @("FAPs", "F2", "F2rl1", "FAPs", 3, 1, 0.6827986666666667, 2.048396, 0.2981510783544719, 0.3185549111733306, 3, 1, 2.662448666666667, 7.987346000000001, 0.9268761553538212, 0.9500328399529846, 1.817916399668444, 16.361247597016, 0.2763491252197888, 0.30263762694297),
  @("FAPs", "F2", "F2rl1", "sCs", 3, 1, 0.6827986666666667, 2.048396, 0.2981510783544719, 0.3185549111733306, 2, 1, 0.210048, 0.420096, 0.07312384464617888, 0.04996716004701549, 0.143420494336, 0.860522966016, 0.02180195313468311, 0.01591728423036062),
  @("M1", "F2", "F2rl1", "FAPs", 1, 0.3333333333333333, 0.073199, 0.219597, 0.03196309812819736, 0.03415047814432848, 3, 1, 2.662448666666667, 7.987346000000001, 0.9268761553538212, 0.9500328399529846, 0.1948885799513333, 1.753997219562, 0.02962583350626049, 0.03244407573720872),
  @("M1", "F2", "F2rl1", "sCs", 1, 0.3333333333333333, 0.073199, 0.219597, 0.03196309812819736, 0.03415047814432848, 2, 1, 0.210048, 0.420096, 0.07312384464617888, 0.04996716004701549, 0.015375303552, 0.09225182131199999, 0.002337264621936875, 0.001706402407119766),
  @("M2", "F2", "F2rl1", "FAPs", 1, 0.3333333333333333, 0.1498013333333333, 0.449404, 0.06541229684924842, 0.0698887574965678, 3, 1, 2.662448666666667, 7.987346000000001, 0.9268761553538212, 0.9500328399529846, 0.3988383601982223, 3.589545241784001, 0.06062909821649425, 0.06639661476524975),
  @("M2", "F2", "F2rl1", "sCs", 1, 0.3333333333333333, 0.1498013333333333, 0.449404, 0.06541229684924842, 0.0698887574965678, 2, 1, 0.210048, 0.420096, 0.07312384464617888, 0.04996716004701549, 0.031465470464, 0.188792822784, 0.004783198632754178, 0.003492142731318057),
  @("sCs", "F2", "F2rl1", "FAPs", 2, 1, 0.440053, 0.8801059999999999, 0.1921536799766067, 0.1368690861791935, 3, 1, 2.662448666666667, 7.987346000000001, 0.9268761553538212, 0.9500328399529846, 1.171618523112667, 7.029711138676, 0.1781026641338057, 0.130030126644589),
  @("sCs", "F2", "F2rl1", "sCs", 2, 1, 0.440053, 0.8801059999999999, 0.1921536799766067, 0.1368690861791935, 2, 1, 0.210048, 0.420096, 0.07312384464617888, 0.04996716004701549, 0.092432252544, 0.369729010176, 0.01405101584280096, 0.00683895953460452)
)

$rowIndex = 2
foreach ($row in $data) {
  $colIndex = 1
  foreach ($val in $row) {
    $ws.Cells.Item($rowIndex, $colIndex).Value = $val
    $colIndex++
  }
  $rowIndex++
}
